$wb = $excel.ActiveWorkbook

$wsPooled = $wb.Worksheets.Item("pooled_effects")
$wsSubgroup = $wb.Worksheets.Item("subgroup_analysis")

# --- pooled_effects: fix "TF-Other" -> "TF-CBT" for the two Raabe, 2022 entries
# in the "excluded" (One ES/study) lists ---

$wsPooled.Range("J3").Value = "Arntz, 2007 (Exposure Therapy vs. Waitlist; PSS-SR); Blanchard, 2003 (TF-CBT vs. Waitlist; CAPS); Brom, 1989 (Exposure Therapy vs. Waitlist; IES); Bryant, 2019 (TF-CBT vs. Waitlist; CAPS); Carlson, 1998 (EMDR vs. TAU; M-PTSD); Carlson, 1998 (EMDR vs. TAU; PSS-SR); Ehlers, 2014 (TF-Cognitive Therapy vs. Waitlist; CAPS); Ehlers, 2014 (TF-Cognitive Therapy vs. Waitlist; CAPS); Foa, 1991 (CBT without trauma focus vs. Waitlist; PSS-I); Foa, 1991 (PE vs. Waitlist; PSS-I); Foa, 1999 (PE vs. Waitlist; PSS-I); Foa, 1999 (CBT without trauma focus vs. Waitlist; PSS-I); Foa, 2005 (TF-CBT vs. Waitlist; PSS-I); Franklin, 2017 (PE vs. TAU; CAPS); McDonagh, 2005 (Other non-trauma-focused vs. Waitlist; CAPS); McGeary, 2022 (TF-Cognitive Therapy vs. TAU; PCL); Power, 2002 (EMDR vs. Waitlist; SI-PTSD); Raabe, 2022 (Other non-trauma-focused vs. Waitlist; CAPS); Raabe, 2022 (TF-CBT vs. Waitlist; CAPS); Raabe, 2022 (TF-CBT vs. Waitlist; PDS); Reger, 2016 (PE vs. Waitlist; CAPS); Resick, 2002 (TF-Cognitive Therapy vs. Other; CAPS); Rothbaum, 2005 (PE vs. Waitlist; CAPS); van den Berg, 2015 (PE vs. Waitlist; CAPS); van Denderen, 2018 (Mixed psychotherapy vs. Waitlist; IES); Wells, 2015 (CBT without trauma focus vs. Waitlist; IES); Wells, 2015 (CBT without trauma focus vs. Waitlist; PDS); Wells, 2015 (PE vs. Waitlist; IES); Zang, 2014 (NET vs. Waitlist; IES)"

$wsPooled.Range("J4").Value = "Arntz, 2007 (Exposure Therapy vs. Waitlist; PSS-SR); Blanchard, 2003 (Other non-trauma-focused vs. Waitlist; CAPS); Brom, 1989 (Other non-trauma-focused vs. Waitlist; IES); Bryant, 2019 (TF-CBT vs. Waitlist; CAPS); Carlson, 1998 (EMDR vs. TAU; IES); Carlson, 1998 (EMDR vs. TAU; M-PTSD); Ehlers, 2014 (TF-Cognitive Therapy vs. Waitlist; CAPS); Ehlers, 2014 (Other non-trauma-focused vs. Waitlist; CAPS); Foa, 1991 (PE vs. Waitlist; PSS-I); Foa, 1991 (Other non-trauma-focused vs. Waitlist; PSS-I); Foa, 1999 (CBT without trauma focus vs. Waitlist; PSS-I); Foa, 1999 (Mixed psychotherapy vs. Waitlist; PSS-I); Foa, 2005 (PE vs. Waitlist; PSS-I); Franklin, 2017 (PE vs. TAU; CAPS); McDonagh, 2005 (TF-CBT vs. Waitlist; CAPS); McGeary, 2022 (CBT without trauma focus vs. TAU; PCL); Power, 2002 (TF-CBT vs. Waitlist; SI-PTSD); Raabe, 2022 (Other non-trauma-focused vs. Waitlist; CAPS); Raabe, 2022 (Other non-trauma-focused vs. Waitlist; PDS); Raabe, 2022 (TF-CBT vs. Waitlist; CAPS); Reger, 2016 (Exposure Therapy vs. Waitlist; CAPS); Resick, 2002 (PE vs. Other; CAPS); Rothbaum, 2005 (EMDR vs. Waitlist; CAPS); van den Berg, 2015 (EMDR vs. Waitlist; CAPS); van Denderen, 2018 (Mixed psychotherapy vs. Waitlist; IES); Wells, 2015 (CBT without trauma focus vs. Waitlist; IES); Wells, 2015 (PE vs. Waitlist; IES); Wells, 2015 (PE vs. Waitlist; PDS); Zang, 2014 (NET vs. Waitlist; IES)"

# --- subgroup_analysis: updated numbers for the condition_arm1.lumped_category
# subgroup rows (TF-CBT / TF-Other / Non-trauma-focused vs waitlist) ---

# Row 7: TF-CBT
$wsSubgroup.Range("D7").Value = 48.0
$wsSubgroup.Range("E7").Value = 1.19
$wsSubgroup.Range("F7").Value = "[0.95; 1.42]"
$wsSubgroup.Range("G7").Value = "'76.9"
$wsSubgroup.Range("H7").Value = "[69.6; 82.4]"
$wsSubgroup.Range("I7").Value = 2.75

# Row 8: TF-Other
$wsSubgroup.Range("D8").Value = 6.0
$wsSubgroup.Range("E8").Value = 1.8
$wsSubgroup.Range("F8").Value = "[0.59; 3.01]"
$wsSubgroup.Range("G8").Value = "'76.1"
$wsSubgroup.Range("H8").Value = "[46.3; 89.3]"
$wsSubgroup.Range("I8").Value = 1.67

# Row 9: Non-trauma-focused
$wsSubgroup.Range("D9").Value = 16.0
$wsSubgroup.Range("E9").Value = 1.14
$wsSubgroup.Range("F9").Value = "[0.68; 1.59]"
$wsSubgroup.Range("G9").Value = "'79.8"
$wsSubgroup.Range("H9").Value = "[67.9; 87.2]"
$wsSubgroup.Range("I9").Value = 2.91
